$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "SIMSmethods" (column C) values for existing rows so that
#     isotope-tagged columns also note the new Ca / CaO methods. ---
# (C2 written first so "d18O10; d13C7; Ca; CaO" lands on the shared-string
#  table before the C12-only variant "d18O10; Ca; CaO".)
$ws.Range("C2").Value  = "d18O10; d13C7; Ca; CaO"
$ws.Range("C3").Value  = "d18O10; d13C7; Ca; CaO"
$ws.Range("C5").Value  = "d18O10; d13C7; Ca; CaO"
$ws.Range("C6").Value  = "d18O10; d13C7; Ca; CaO"
$ws.Range("C8").Value  = "d18O10; d13C7; Ca; CaO"
$ws.Range("C10").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C12").Value = "d18O10; Ca; CaO"
$ws.Range("C13").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C14").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C15").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C16").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C17").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C18").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("C19").Value = "d18O10; d13C7; Ca; CaO"

# --- Append new column-dictionary rows for the Ca/CaO isotope columns and
#     the is_standard flag. ---

# Row 40: d44Cameas
$ws.Range("A40").Value = "d44Cameas"
$ws.Range("C40").Value = "Ca; CaO"
$ws.Range("B40").Value = "\u03B444Ca \u2030 measured"
$ws.Range("D40").Value = "permille"
$ws.Range("E40").Value = "Numeric"
$ws.Range("C40").Font.Color = 0
$ws.Range("E40").Font.Color = 0

# Row 41: Ca40Cps
$ws.Range("B41").Value = "40Ca (Gcps)"
$ws.Range("C41").Value = "Ca"
$ws.Range("A41").Value = "Ca40Cps"
$ws.Range("D41").Value = "Gcps"
$ws.Range("E41").Value = "Numeric"
$ws.Range("E41").Font.Color = 0

# Row 42: is_standard
$ws.Range("A42").Value = "is_standard"
$ws.Range("D42").Value = "Boolean"
$ws.Range("C42").Value = "d18O10; d13C7; Ca; CaO"
$ws.Range("E42").Value = "Text"
$ws.Range("E42").Font.Color = 0

# --- Restore the view state that Excel persisted after the edit. ---
$ws.Range("E45").Select()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 2
